# Scheduled runner update: refresh market-price-derived profit figures
# across the per-job Leve sheets (currentAveragePrice*, LevePrice*,
# LeveProfit* columns H:N) with newly pulled values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 21778.166
$ws.Range("J21").Value = 13144.444
$ws.Range("L21").Value = 13144.444
$ws.Range("N21").Value = -14080.444
$ws.Range("H23").Value = 21778.166
$ws.Range("J23").Value = 13144.444
$ws.Range("L23").Value = 13144.444
$ws.Range("N23").Value = -13612.444
$ws.Range("H28").Value = 1143.9445
$ws.Range("I28").Value = 349.16666
$ws.Range("J28").Value = 2733.5
$ws.Range("K28").Value = 349.16666
$ws.Range("L28").Value = 2733.5
$ws.Range("M28").Value = 135.83334
$ws.Range("N28").Value = -3703.5
$ws.Range("H111").Value = 333333700
$ws.Range("I111").Value = 500000260
$ws.Range("J111").Value = 600
$ws.Range("K111").Value = 1500000780
$ws.Range("L111").Value = 1800
$ws.Range("M111").Value = -1499997713
$ws.Range("N111").Value = -7934
$ws.Range("H132").Value = 2292.5442
$ws.Range("I132").Value = 2012.0546
$ws.Range("J132").Value = 3479.2307
$ws.Range("K132").Value = 6036.1638
$ws.Range("L132").Value = 10437.6921
$ws.Range("M132").Value = -3506.1638
$ws.Range("N132").Value = -15497.6921

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20837082
$ws.Range("I32").Value = 22730112
$ws.Range("J32").Value = 13753.5
$ws.Range("K32").Value = 22730112
$ws.Range("L32").Value = 13753.5
$ws.Range("M32").Value = -22729825
$ws.Range("N32").Value = -14327.5
$ws.Range("H45").Value = 1749.3667
$ws.Range("I45").Value = 930
$ws.Range("J45").Value = 1807.8928
$ws.Range("K45").Value = 930
$ws.Range("L45").Value = 1807.8928
$ws.Range("M45").Value = -553
$ws.Range("N45").Value = -2561.8928
$ws.Range("H61").Value = 2319.0667
$ws.Range("I61").Value = 1351
$ws.Range("K61").Value = 1351
$ws.Range("M61").Value = -1139
$ws.Range("H63").Value = 2680
$ws.Range("I63").Value = 2270
$ws.Range("J63").Value = 3500
$ws.Range("K63").Value = 2270
$ws.Range("L63").Value = 3500
$ws.Range("M63").Value = -1584
$ws.Range("N63").Value = -4872
$ws.Range("H66").Value = 2680
$ws.Range("I66").Value = 2270
$ws.Range("J66").Value = 3500
$ws.Range("K66").Value = 11350
$ws.Range("L66").Value = 17500
$ws.Range("M66").Value = -7918
$ws.Range("N66").Value = -24364
$ws.Range("H132").Value = 1967.4642
$ws.Range("I132").Value = 1113.579
$ws.Range("J132").Value = 3770.111
$ws.Range("K132").Value = 3340.737
$ws.Range("L132").Value = 11310.333
$ws.Range("M132").Value = -810.7370000000001
$ws.Range("N132").Value = -16370.333
$ws.Range("H136").Value = 2319.0667
$ws.Range("I136").Value = 1351
$ws.Range("K136").Value = 4053
$ws.Range("M136").Value = -1503

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1411.4736
$ws.Range("I134").Value = 1407.5652
$ws.Range("J134").Value = 1427.8182
$ws.Range("K134").Value = 4222.6956
$ws.Range("L134").Value = 4283.4546
$ws.Range("M134").Value = -1687.6956
$ws.Range("N134").Value = -9353.454600000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1854.4546
$ws.Range("I132").Value = 1435.8928
$ws.Range("J132").Value = 4198.4
$ws.Range("K132").Value = 4307.678400000001
$ws.Range("L132").Value = 12595.2
$ws.Range("M132").Value = -1777.678400000001
$ws.Range("N132").Value = -17655.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 5.111111
$ws.Range("J12").Value = 3
$ws.Range("L12").Value = 9
$ws.Range("N12").Value = -355
$ws.Range("H74").Value = 3810.8
$ws.Range("J74").Value = 8007.5
$ws.Range("L74").Value = 24022.5
$ws.Range("N74").Value = -26144.5
$ws.Range("H77").Value = 3810.8
$ws.Range("J77").Value = 8007.5
$ws.Range("L77").Value = 72067.5
$ws.Range("N77").Value = -82675.5
$ws.Range("H113").Value = 585080.1
$ws.Range("I113").Value = 1014817.25
$ws.Range("J113").Value = 637.6
$ws.Range("K113").Value = 3044451.75
$ws.Range("L113").Value = 1912.8
$ws.Range("M113").Value = -3042281.75
$ws.Range("N113").Value = -6252.8
$ws.Range("H131").Value = 768.15
$ws.Range("I131").Value = 460
$ws.Range("J131").Value = 900.2143
$ws.Range("K131").Value = 1380
$ws.Range("L131").Value = 2700.6429
$ws.Range("M131").Value = 3660
$ws.Range("N131").Value = -12780.6429

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3510.3635
$ws.Range("I102").Value = 3725
$ws.Range("J102").Value = 2938
$ws.Range("K102").Value = 3725
$ws.Range("L102").Value = 2938
$ws.Range("M102").Value = -2103
$ws.Range("N102").Value = -6182
$ws.Range("H122").Value = 16668551
$ws.Range("I122").Value = 50001250
$ws.Range("J122").Value = 2202
$ws.Range("K122").Value = 150003750
$ws.Range("L122").Value = 6606
$ws.Range("M122").Value = -150001300
$ws.Range("N122").Value = -11506
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H132").Value = 2007.8474
$ws.Range("I132").Value = 1736.8541
$ws.Range("K132").Value = 5210.5623
$ws.Range("M132").Value = -2680.5623
$ws.Range("H139").Value = 30000
$ws.Range("J139").Value = 30000
$ws.Range("L139").Value = 30000
$ws.Range("N139").Value = -40280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5000
$ws.Range("I7").Value = 5000
$ws.Range("J7").Value = 5000
$ws.Range("K7").Value = 5000
$ws.Range("L7").Value = 5000
$ws.Range("M7").Value = -4888
$ws.Range("N7").Value = -5224
$ws.Range("H36").Value = 72500
$ws.Range("J36").Value = 72500
$ws.Range("L36").Value = 72500
$ws.Range("N36").Value = -73624
$ws.Range("H40").Value = 4809.85
$ws.Range("I40").Value = 4019.4666
$ws.Range("J40").Value = 7181
$ws.Range("K40").Value = 4019.4666
$ws.Range("L40").Value = 7181
$ws.Range("M40").Value = -3883.4666
$ws.Range("N40").Value = -7453
$ws.Range("H122").Value = 3403.75
$ws.Range("I122").Value = 3250
$ws.Range("K122").Value = 9750
$ws.Range("M122").Value = -7300
$ws.Range("H126").Value = 5000
$ws.Range("I126").Value = 5000
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 15000
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -12530
$ws.Range("N126").Value = -19940

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2667.4443
$ws.Range("I122").Value = 2545.6924
$ws.Range("K122").Value = 7637.0772
$ws.Range("M122").Value = -5187.0772
$ws.Range("H132").Value = 1164.3392
$ws.Range("I132").Value = 731.9778
$ws.Range("J132").Value = 2933.0908
$ws.Range("K132").Value = 2195.9334
$ws.Range("L132").Value = 8799.2724
$ws.Range("M132").Value = 334.0666000000001
$ws.Range("N132").Value = -13859.2724
$ws.Range("H138").Value = 28329.334
$ws.Range("J138").Value = 28329.334
$ws.Range("L138").Value = 28329.334
$ws.Range("N138").Value = -38609.334
